$d = $word.ActiveDocument

# Each of these Find/Replace calls targets text that already spans the
# seam between two adjacent runs (a text run immediately followed by a
# single-space run, or vice-versa). Replacing the matched range with the
# identical text causes the runtime to re-materialize it as one run,
# merging the previously-split runs without altering any visible text,
# comment ranges, or comment references.

# Paragraph 1 (FirstParagraph): "I want" + " " -> "I want "
$d.Content.Find.Execute("I want ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "I want ", 2)

# Paragraph 1 (FirstParagraph): "some text to have a comment" + " " -> "some text to have a comment "
$d.Content.Find.Execute("some text to have a comment ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "some text to have a comment ", 2)

# Paragraph 2 (BodyText): "This is" + " " -> "This is "
$d.Content.Find.Execute("This is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This is ", 2)

# Paragraph 3 (BodyText): " " + "is this." -> " is this."
$d.Content.Find.Execute(" is this.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " is this.", 2)

# Paragraph 4 (BodyText): "One" + " " -> "One "
$d.Content.Find.Execute("One ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "One ", 2)

# Paragraph 4 (BodyText): ". And this is one with a" + " " -> ". And this is one with a "
$d.Content.Find.Execute(". And this is one with a ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". And this is one with a ", 2)
